$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $ws.Range("A1").Style
}

$ws.Range("D2").Value = "27.102.00"
$ws.Range("E2").Value = "  -2.61%  "

$ws.Range("D3").Value = "1.867.52"
$ws.Range("E3").Value = "  -2.07%  "

Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue $ws.Range("D5") "306.09"
$ws.Range("E5").Value = "  -2.19%  "

Set-TextValue $ws.Range("D6") "0.9998"
$ws.Range("E6").Value = "  -0.07%  "

Set-TextValue $ws.Range("D7") "0.5166"
$ws.Range("E7").Value = "  -1.31%  "

Set-TextValue $ws.Range("D8") "0.3771"
$ws.Range("E8").Value = "  -0.40%  "

Set-TextValue $ws.Range("D9") "0.07163"

Set-TextValue $ws.Range("D10") "0.8905"
$ws.Range("E10").Value = "  -2.06%  "

Set-TextValue $ws.Range("D11") "20.72"
$ws.Range("E11").Value = "  -2.62%  "

Set-TextValue $ws.Range("D12") "0.07606"
$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "1.872.07"
$ws.Range("E13").Value = "  -2.37%  "

Set-TextValue $ws.Range("D14") "5.315"
$ws.Range("E14").Value = "  -2.56%  "

Set-TextValue $ws.Range("D15") "89.81"
$ws.Range("E15").Value = "  -2.54%  "

Set-TextValue $ws.Range("D16") "1.001"
$ws.Range("E16").Value = "  +0.04%  "

Set-TextValue $ws.Range("D17") "0.000008497"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("E18").Value = "  -3.24%  "

Set-TextValue $ws.Range("D19") "0.9997"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").Value = "27.117.17"
$ws.Range("E20").Value = "  -2.64%  "

Set-TextValue $ws.Range("D21") "5.032"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").Value = "2.119.89"
$ws.Range("E22").Value = "  -2.58%  "

$ws.Range("E23").Value = "  -3.28%  "

Set-TextValue $ws.Range("D24") "6.470"
$ws.Range("E24").Value = "  -2.63%  "

Set-TextValue $ws.Range("D25") "1.839"
$ws.Range("E25").Value = "  -1.83%  "

Set-TextValue $ws.Range("D26") "147.56"
$ws.Range("E26").Value = "  -4.08%  "

Set-TextValue $ws.Range("D27") "17.96"
$ws.Range("E27").Value = "  -2.01%  "

Set-TextValue $ws.Range("D28") "2.095"
$ws.Range("E28").Value = "  -3.44%  "

Set-TextValue $ws.Range("D29") "112.90"
$ws.Range("E29").Value = "  -1.65%  "

Set-TextValue $ws.Range("D30") "4.664"
$ws.Range("E30").Value = "  -4.10%  "

Set-TextValue $ws.Range("D31") "4.675"
$ws.Range("E31").Value = "  -3.83%  "

Set-TextValue $ws.Range("D32") "0.09149"
$ws.Range("E32").Value = "  +1.35%  "

Set-TextValue $ws.Range("D33") "0.05114"
$ws.Range("E33").Value = "  -3.27%  "

$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("E35").Value = "  -6.36%  "

Set-TextValue $ws.Range("D36") "0.7266"
$ws.Range("E36").Value = "  -6.93%  "

Set-TextValue $ws.Range("D37") "0.02035"
$ws.Range("E37").Value = "  -2.82%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "3.070"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "2.501"
$ws.Range("E39").Value = "  -4.05%  "

$ws.Range("E40").Value = "  -1.92%  "

Set-TextValue $ws.Range("D41") "0.5327"
$ws.Range("E41").Value = "  -4.78%  "

Set-TextValue $ws.Range("D42") "6.491"
$ws.Range("E42").Value = "  -3.54%  "

Set-TextValue $ws.Range("D43") "116.42"
$ws.Range("E43").Value = "  +0.81%  "

Set-TextValue $ws.Range("D44") "8.294"
$ws.Range("E44").Value = "  -3.26%  "

Set-TextValue $ws.Range("D45") "0.1467"
$ws.Range("E45").Value = "  -3.11%  "

Set-TextValue $ws.Range("D46") "0.4633"
$ws.Range("E46").Value = "  -3.60%  "

Set-TextValue $ws.Range("D47") "0.9991"
$ws.Range("E47").Value = "  -0.15%  "

Set-TextValue $ws.Range("D48") "9.984"
$ws.Range("E48").Value = "  -5.07%  "

Set-TextValue $ws.Range("D49") "1.572"
$ws.Range("E49").Value = "  -3.15%  "

Set-TextValue $ws.Range("D50") "36.58"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("E51").Value = "  -4.98%  "
